$d = $word.ActiveDocument

# The test bench's expected simulation duration changed from ~5000 ns to
# ~9500 ns. That value lives inside an OMath (equation) run, which Find/
# Replace cannot reach, so rebuild the equation's contents via InsertXML.

# Word keeps a single "_GoBack" bookmark marking the last edit location;
# remove the old one (after "...answer packets.") so it can be
# re-created at the new edit site inside the equation.
try {
    $d.Bookmarks.Item("_GoBack").Delete()
} catch {
}

$om = $d.OMaths.Item(2)
$r = $om.Range
$xml = '<m:oMath>' +
  '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Courier New"/></w:rPr><m:t>' + [char]0x2248 + '</m:t></m:r>' +
  '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Courier New"/></w:rPr><m:t>95</m:t></m:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Courier New"/></w:rPr><m:t>00 ns</m:t></m:r>' +
  '</m:oMath>'
$r.InsertXML($xml)
